$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI pipeline was re-run against the new TPM matrix. For the Wnt2->Fzd7
# pair this produced a new "ECs" sending-cluster block (one row per target
# cluster: ECs, FAPs, MuSCs) ahead of the pre-existing "FAPs" sending-cluster
# block, whose derived-specificity figures also shifted because they are
# computed relative to the whole (now different) expression matrix.

# Make room: push the current 3 data rows (FAPs block) down from 2-4 to 5-7.
$ws.Rows("2:4").Insert()
# Row-insert in this engine copies the formatting of the row above (the bold
# header row); strip that back off so the new rows use the plain data style
# that every other data row already has.
$ws.Range("A2:T4").ClearFormats()

# New rows: Sending cluster = ECs
$ws.Cells.Item(2, 1).Value  = "ECs"
$ws.Cells.Item(2, 2).Value  = "Wnt2"
$ws.Cells.Item(2, 3).Value  = "Fzd7"
$ws.Cells.Item(2, 4).Value  = "ECs"
$ws.Cells.Item(2, 5).Value  = 1
$ws.Cells.Item(2, 6).Value  = 0.3333333333333333
$ws.Cells.Item(2, 7).Value  = 0.01070233333333333
$ws.Cells.Item(2, 8).Value  = 0.032107
$ws.Cells.Item(2, 9).Value  = 0.004227647500550067
$ws.Cells.Item(2, 10).Value = 0.004227647500550067
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.6068319999999999
$ws.Cells.Item(2, 14).Value = 1.820496
$ws.Cells.Item(2, 15).Value = 0.03392274820144286
$ws.Cells.Item(2, 16).Value = 0.03392274820144286
$ws.Cells.Item(2, 17).Value = 0.006494518341333333
$ws.Cells.Item(2, 18).Value = 0.05845066507199999
$ws.Cells.Item(2, 19).Value = 0.0001434134216456192
$ws.Cells.Item(2, 20).Value = 0.0001434134216456192

$ws.Cells.Item(3, 1).Value  = "ECs"
$ws.Cells.Item(3, 2).Value  = "Wnt2"
$ws.Cells.Item(3, 3).Value  = "Fzd7"
$ws.Cells.Item(3, 4).Value  = "FAPs"
$ws.Cells.Item(3, 5).Value  = 1
$ws.Cells.Item(3, 6).Value  = 0.3333333333333333
$ws.Cells.Item(3, 7).Value  = 0.01070233333333333
$ws.Cells.Item(3, 8).Value  = 0.032107
$ws.Cells.Item(3, 9).Value  = 0.004227647500550067
$ws.Cells.Item(3, 10).Value = 0.004227647500550067
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 8.058662
$ws.Cells.Item(3, 14).Value = 24.175986
$ws.Cells.Item(3, 15).Value = 0.4504903529585388
$ws.Cells.Item(3, 16).Value = 0.4504903529585388
$ws.Cells.Item(3, 17).Value = 0.08624648694466666
$ws.Cells.Item(3, 18).Value = 0.7762183825019999
$ws.Cells.Item(3, 19).Value = 0.001904514414707084
$ws.Cells.Item(3, 20).Value = 0.001904514414707084

$ws.Cells.Item(4, 1).Value  = "ECs"
$ws.Cells.Item(4, 2).Value  = "Wnt2"
$ws.Cells.Item(4, 3).Value  = "Fzd7"
$ws.Cells.Item(4, 4).Value  = "MuSCs"
$ws.Cells.Item(4, 5).Value  = 1
$ws.Cells.Item(4, 6).Value  = 0.3333333333333333
$ws.Cells.Item(4, 7).Value  = 0.01070233333333333
$ws.Cells.Item(4, 8).Value  = 0.032107
$ws.Cells.Item(4, 9).Value  = 0.004227647500550067
$ws.Cells.Item(4, 10).Value = 0.004227647500550067
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 9.223151
$ws.Cells.Item(4, 14).Value = 27.669453
$ws.Cells.Item(4, 15).Value = 0.5155868988400183
$ws.Cells.Item(4, 16).Value = 0.5155868988400183
$ws.Cells.Item(4, 17).Value = 0.09870923638566666
$ws.Cells.Item(4, 18).Value = 0.8883831274709999
$ws.Cells.Item(4, 19).Value = 0.002179719664197363
$ws.Cells.Item(4, 20).Value = 0.002179719664197363

# Existing rows (now 5-7, sending cluster = FAPs): ligand/receptor-detection
# columns (E-J) are unchanged by the TPM update, but the derived-specificity
# and edge-weight columns (M-T), which depend on the full expression matrix,
# were recomputed.
$ws.Cells.Item(5, 1).Value  = "FAPs"
$ws.Cells.Item(5, 2).Value  = "Wnt2"
$ws.Cells.Item(5, 3).Value  = "Fzd7"
$ws.Cells.Item(5, 4).Value  = "ECs"
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 6).Value  = 1
$ws.Cells.Item(5, 7).Value  = 2.520808
$ws.Cells.Item(5, 8).Value  = 7.562424
$ws.Cells.Item(5, 9).Value  = 0.99577235249945
$ws.Cells.Item(5, 10).Value = 0.99577235249945
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6068319999999999
$ws.Cells.Item(5, 14).Value = 1.820496
$ws.Cells.Item(5, 15).Value = 0.03392274820144286
$ws.Cells.Item(5, 16).Value = 0.03392274820144286
$ws.Cells.Item(5, 17).Value = 1.529706960256
$ws.Cells.Item(5, 18).Value = 13.767362642304
$ws.Cells.Item(5, 19).Value = 0.03377933477979725
$ws.Cells.Item(5, 20).Value = 0.03377933477979725

$ws.Cells.Item(6, 1).Value  = "FAPs"
$ws.Cells.Item(6, 2).Value  = "Wnt2"
$ws.Cells.Item(6, 3).Value  = "Fzd7"
$ws.Cells.Item(6, 4).Value  = "FAPs"
$ws.Cells.Item(6, 5).Value  = 3
$ws.Cells.Item(6, 6).Value  = 1
$ws.Cells.Item(6, 7).Value  = 2.520808
$ws.Cells.Item(6, 8).Value  = 7.562424
$ws.Cells.Item(6, 9).Value  = 0.99577235249945
$ws.Cells.Item(6, 10).Value = 0.99577235249945
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 8.058662
$ws.Cells.Item(6, 14).Value = 24.175986
$ws.Cells.Item(6, 15).Value = 0.4504903529585388
$ws.Cells.Item(6, 16).Value = 0.4504903529585388
$ws.Cells.Item(6, 17).Value = 20.314339638896
$ws.Cells.Item(6, 18).Value = 182.829056750064
$ws.Cells.Item(6, 19).Value = 0.4485858385438317
$ws.Cells.Item(6, 20).Value = 0.4485858385438317

$ws.Cells.Item(7, 1).Value  = "FAPs"
$ws.Cells.Item(7, 2).Value  = "Wnt2"
$ws.Cells.Item(7, 3).Value  = "Fzd7"
$ws.Cells.Item(7, 4).Value  = "MuSCs"
$ws.Cells.Item(7, 5).Value  = 3
$ws.Cells.Item(7, 6).Value  = 1
$ws.Cells.Item(7, 7).Value  = 2.520808
$ws.Cells.Item(7, 8).Value  = 7.562424
$ws.Cells.Item(7, 9).Value  = 0.99577235249945
$ws.Cells.Item(7, 10).Value = 0.99577235249945
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.223151
$ws.Cells.Item(7, 14).Value = 27.669453
$ws.Cells.Item(7, 15).Value = 0.5155868988400183
$ws.Cells.Item(7, 16).Value = 0.5155868988400183
$ws.Cells.Item(7, 17).Value = 23.249792826008
$ws.Cells.Item(7, 18).Value = 209.248135434072
$ws.Cells.Item(7, 19).Value = 0.5134071791758209
$ws.Cells.Item(7, 20).Value = 0.5134071791758209
